$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 596 (shifts existing rows 596-640 down to 597-641)
$ws.Rows.Item(596).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A596").Value = 9
$ws.Range("B596").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C596").Value = "Metropolitana"
$ws.Range("D596").Value = 45223
$ws.Range("E596").Value = 13
$ws.Range("F596").Value = 100112052
$ws.Range("G596").Value = "Albahaca"
$ws.Range("H596").Value = "Sin especificar"
$ws.Range("I596").Value = "Primera"
$ws.Range("J596").Value = 160
$ws.Range("K596").Value = 5000
$ws.Range("L596").Value = 5000
$ws.Range("M596").Value = 5000
$ws.Range("N596").Value = "$/docena de matas"
$ws.Range("O596").Value = "Provincia de Chacabuco"
$ws.Range("P596").Value = 833
$ws.Range("Q596").Value = 6
$ws.Range("R596").Value = "Hortaliza"
